# Removed priority on testcase level and updated tesng.xml
#
# The "execute" column (C) previously had several setup/pre-req rows
# (login, createsJobTitle, createsPayGrade, createsEmploymentStatus,
# createsJobCategories, updatesGeneralInformationData) marked "no" so they
# were skipped at run time. They are flipped to "yes" so those test cases
# run; this also drops "no" from the shared-string table since nothing else
# references it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "yes"
$ws.Range("C4").Value = "yes"
$ws.Range("C6").Value = "yes"
$ws.Range("C8").Value = "yes"
$ws.Range("C10").Value = "yes"
$ws.Range("C12").Value = "yes"

# Clear the saved cursor/selection (the sheet no longer needs to reopen
# focused on H15); reset it back to the top-left cell.
$ws.Range("A1").Select()

$wb.Save()
